# Update LR-pairs TPM-derived statistics (columns G-J, M-P, Q-T) for rows 2-26
# with recomputed values from the new TPM input, per the commit "update scripts wuth new tpm".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 22.95531766666667
$ws.Range("H2").Value = 68.865953
$ws.Range("I2").Value = 0.1720020945576478
$ws.Range("J2").Value = 0.1720020945576478
$ws.Range("M2").Value = 1.620350333333333
$ws.Range("N2").Value = 4.861051
$ws.Range("O2").Value = 0.0725197794467048
$ws.Range("P2").Value = 0.07251977944670479
$ws.Range("Q2").Value = 37.1956566329559
$ws.Range("R2").Value = 334.760909696603
$ws.Range("S2").Value = 0.01247355396169189
$ws.Range("T2").Value = 0.01247355396169188

# Row 3
$ws.Range("G3").Value = 22.95531766666667
$ws.Range("H3").Value = 68.865953
$ws.Range("I3").Value = 0.1720020945576478
$ws.Range("J3").Value = 0.1720020945576478
$ws.Range("M3").Value = 2.170845333333333
$ws.Range("N3").Value = 6.512536
$ws.Range("O3").Value = 0.09715752300453648
$ws.Range("P3").Value = 0.09715752300453648
$ws.Range("Q3").Value = 49.83244423186756
$ws.Range("R3").Value = 448.491998086808
$ws.Range("S3").Value = 0.01671129745881313
$ws.Range("T3").Value = 0.01671129745881313

# Row 4
$ws.Range("G4").Value = 22.95531766666667
$ws.Range("H4").Value = 68.865953
$ws.Range("I4").Value = 0.1720020945576478
$ws.Range("J4").Value = 0.1720020945576478
$ws.Range("M4").Value = 13.65106133333333
$ws.Range("N4").Value = 40.953184
$ws.Range("O4").Value = 0.6109616770777183
$ws.Range("P4").Value = 0.6109616770777183
$ws.Range("Q4").Value = 313.3644493938169
$ws.Range("R4").Value = 2820.280044544352
$ws.Range("S4").Value = 0.1050866881518208
$ws.Range("T4").Value = 0.1050866881518208

# Row 5
$ws.Range("G5").Value = 22.95531766666667
$ws.Range("H5").Value = 68.865953
$ws.Range("I5").Value = 0.1720020945576478
$ws.Range("J5").Value = 0.1720020945576478
$ws.Range("M5").Value = 0.6533493333333333
$ws.Range("N5").Value = 1.960048
$ws.Range("O5").Value = 0.02924105274043717
$ws.Range("P5").Value = 0.02924105274043717
$ws.Range("Q5").Value = 14.99784149397156
$ws.Range("R5").Value = 134.980573445744
$ws.Range("S5").Value = 0.005029522318425842
$ws.Range("T5").Value = 0.005029522318425842

# Row 6
$ws.Range("G6").Value = 22.95531766666667
$ws.Range("H6").Value = 68.865953
$ws.Range("I6").Value = 0.1720020945576478
$ws.Range("J6").Value = 0.1720020945576478
$ws.Range("M6").Value = 4.247957666666667
$ws.Range("N6").Value = 12.743873
$ws.Range("O6").Value = 0.1901199677306032
$ws.Range("P6").Value = 0.1901199677306032
$ws.Range("Q6").Value = 97.51321767288547
$ws.Range("R6").Value = 877.6189590559691
$ws.Range("S6").Value = 0.03270103266689617
$ws.Range("T6").Value = 0.03270103266689616

# Row 7
$ws.Range("G7").Value = 62.21411366666666
$ws.Range("H7").Value = 186.642341
$ws.Range("I7").Value = 0.4661646602805707
$ws.Range("J7").Value = 0.4661646602805707
$ws.Range("M7").Value = 1.620350333333333
$ws.Range("N7").Value = 4.861051
$ws.Range("O7").Value = 0.0725197794467048
$ws.Range("P7").Value = 0.07251977944670479
$ws.Range("Q7").Value = 100.8086598178212
$ws.Range("R7").Value = 907.2779383603909
$ws.Range("S7").Value = 0.03380615834939506
$ws.Range("T7").Value = 0.03380615834939505

# Row 8
$ws.Range("G8").Value = 62.21411366666666
$ws.Range("H8").Value = 186.642341
$ws.Range("I8").Value = 0.4661646602805707
$ws.Range("J8").Value = 0.4661646602805707
$ws.Range("M8").Value = 2.170845333333333
$ws.Range("N8").Value = 6.512536
$ws.Range("O8").Value = 0.09715752300453648
$ws.Range("P8").Value = 0.09715752300453648
$ws.Range("Q8").Value = 135.0572183207529
$ws.Range("R8").Value = 1215.514964886776
$ws.Range("S8").Value = 0.04529140370511148
$ws.Range("T8").Value = 0.04529140370511148

# Row 9
$ws.Range("G9").Value = 62.21411366666666
$ws.Range("H9").Value = 186.642341
$ws.Range("I9").Value = 0.4661646602805707
$ws.Range("J9").Value = 0.4661646602805707
$ws.Range("M9").Value = 13.65106133333333
$ws.Range("N9").Value = 40.953184
$ws.Range("O9").Value = 0.6109616770777183
$ws.Range("P9").Value = 0.6109616770777183
$ws.Range("Q9").Value = 849.2886814626381
$ws.Range("R9").Value = 7643.598133163743
$ws.Range("S9").Value = 0.2848087426393823
$ws.Range("T9").Value = 0.2848087426393823

# Row 10
$ws.Range("G10").Value = 62.21411366666666
$ws.Range("H10").Value = 186.642341
$ws.Range("I10").Value = 0.4661646602805707
$ws.Range("J10").Value = 0.4661646602805707
$ws.Range("M10").Value = 0.6533493333333333
$ws.Range("N10").Value = 1.960048
$ws.Range("O10").Value = 0.02924105274043717
$ws.Range("P10").Value = 0.02924105274043717
$ws.Range("Q10").Value = 40.64754968804089
$ws.Range("R10").Value = 365.827947192368
$ws.Range("S10").Value = 0.01363114541699214
$ws.Range("T10").Value = 0.01363114541699214

# Row 11
$ws.Range("G11").Value = 62.21411366666666
$ws.Range("H11").Value = 186.642341
$ws.Range("I11").Value = 0.4661646602805707
$ws.Range("J11").Value = 0.4661646602805707
$ws.Range("M11").Value = 4.247957666666667
$ws.Range("N11").Value = 12.743873
$ws.Range("O11").Value = 0.1901199677306032
$ws.Range("P11").Value = 0.1901199677306032
$ws.Range("Q11").Value = 264.2829211251881
$ws.Range("R11").Value = 2378.546290126693
$ws.Range("S11").Value = 0.08862721016968969
$ws.Range("T11").Value = 0.08862721016968969

# Row 12
$ws.Range("G12").Value = 7.783044333333334
$ws.Range("H12").Value = 23.349133
$ws.Range("I12").Value = 0.05831763893698088
$ws.Range("J12").Value = 0.05831763893698089
$ws.Range("M12").Value = 1.620350333333333
$ws.Range("N12").Value = 4.861051
$ws.Range("O12").Value = 0.0725197794467048
$ws.Range("P12").Value = 0.07251977944670479
$ws.Range("Q12").Value = 12.61125847986478
$ws.Range("R12").Value = 113.501326318783
$ws.Range("S12").Value = 0.004229182313562418
$ws.Range("T12").Value = 0.004229182313562418

# Row 13
$ws.Range("G13").Value = 7.783044333333334
$ws.Range("H13").Value = 23.349133
$ws.Range("I13").Value = 0.05831763893698088
$ws.Range("J13").Value = 0.05831763893698089
$ws.Range("M13").Value = 2.170845333333333
$ws.Range("N13").Value = 6.512536
$ws.Range("O13").Value = 0.09715752300453648
$ws.Range("P13").Value = 0.09715752300453648
$ws.Range("Q13").Value = 16.89578547014311
$ws.Range("R13").Value = 152.062069231288
$ws.Range("S13").Value = 0.005665997346589972
$ws.Range("T13").Value = 0.005665997346589973

# Row 14
$ws.Range("G14").Value = 7.783044333333334
$ws.Range("H14").Value = 23.349133
$ws.Range("I14").Value = 0.05831763893698088
$ws.Range("J14").Value = 0.05831763893698089
$ws.Range("M14").Value = 13.65106133333333
$ws.Range("N14").Value = 40.953184
$ws.Range("O14").Value = 0.6109616770777183
$ws.Range("P14").Value = 0.6109616770777183
$ws.Range("Q14").Value = 106.2468155543858
$ws.Range("R14").Value = 956.2213399894721
$ws.Range("S14").Value = 0.03562984248815068
$ws.Range("T14").Value = 0.03562984248815069

# Row 15
$ws.Range("G15").Value = 7.783044333333334
$ws.Range("H15").Value = 23.349133
$ws.Range("I15").Value = 0.05831763893698088
$ws.Range("J15").Value = 0.05831763893698089
$ws.Range("M15").Value = 0.6533493333333333
$ws.Range("N15").Value = 1.960048
$ws.Range("O15").Value = 0.02924105274043717
$ws.Range("P15").Value = 0.02924105274043717
$ws.Range("Q15").Value = 5.085046826487111
$ws.Range("R15").Value = 45.765421438384
$ws.Range("S15").Value = 0.00170526915585403
$ws.Range("T15").Value = 0.00170526915585403

# Row 16
$ws.Range("G16").Value = 7.783044333333334
$ws.Range("H16").Value = 23.349133
$ws.Range("I16").Value = 0.05831763893698088
$ws.Range("J16").Value = 0.05831763893698089
$ws.Range("M16").Value = 4.247957666666667
$ws.Range("N16").Value = 12.743873
$ws.Range("O16").Value = 0.1901199677306032
$ws.Range("P16").Value = 0.1901199677306032
$ws.Range("Q16").Value = 33.0620428457899
$ws.Range("R16").Value = 297.558385612109
$ws.Range("S16").Value = 0.01108734763282377
$ws.Range("T16").Value = 0.01108734763282377

# Row 17
$ws.Range("G17").Value = 30.44016466666666
$ws.Range("H17").Value = 91.320494
$ws.Range("I17").Value = 0.2280853681650076
$ws.Range("J17").Value = 0.2280853681650076
$ws.Range("M17").Value = 1.620350333333333
$ws.Range("N17").Value = 4.861051
$ws.Range("O17").Value = 0.0725197794467048
$ws.Range("P17").Value = 0.07251977944670479
$ws.Range("Q17").Value = 49.32373096435489
$ws.Range("R17").Value = 443.913578679194
$ws.Range("S17").Value = 0.01654070059434682
$ws.Range("T17").Value = 0.01654070059434681

# Row 18
$ws.Range("G18").Value = 30.44016466666666
$ws.Range("H18").Value = 91.320494
$ws.Range("I18").Value = 0.2280853681650076
$ws.Range("J18").Value = 0.2280853681650076
$ws.Range("M18").Value = 2.170845333333333
$ws.Range("N18").Value = 6.512536
$ws.Range("O18").Value = 0.09715752300453648
$ws.Range("P18").Value = 0.09715752300453648
$ws.Range("Q18").Value = 66.08088941253155
$ws.Range("R18").Value = 594.728004712784
$ws.Range("S18").Value = 0.0221602094044899
$ws.Range("T18").Value = 0.0221602094044899

# Row 19
$ws.Range("G19").Value = 30.44016466666666
$ws.Range("H19").Value = 91.320494
$ws.Range("I19").Value = 0.2280853681650076
$ws.Range("J19").Value = 0.2280853681650076
$ws.Range("M19").Value = 13.65106133333333
$ws.Range("N19").Value = 40.953184
$ws.Range("O19").Value = 0.6109616770777183
$ws.Range("P19").Value = 0.6109616770777183
$ws.Range("Q19").Value = 415.5405548614328
$ws.Range("R19").Value = 3739.864993752896
$ws.Range("S19").Value = 0.1393514190509819
$ws.Range("T19").Value = 0.1393514190509819

# Row 20
$ws.Range("G20").Value = 30.44016466666666
$ws.Range("H20").Value = 91.320494
$ws.Range("I20").Value = 0.2280853681650076
$ws.Range("J20").Value = 0.2280853681650076
$ws.Range("M20").Value = 0.6533493333333333
$ws.Range("N20").Value = 1.960048
$ws.Range("O20").Value = 0.02924105274043717
$ws.Range("P20").Value = 0.02924105274043717
$ws.Range("Q20").Value = 19.88806129152355
$ws.Range("R20").Value = 178.992551623712
$ws.Range("S20").Value = 0.006669456279835016
$ws.Range("T20").Value = 0.006669456279835016

# Row 21
$ws.Range("G21").Value = 30.44016466666666
$ws.Range("H21").Value = 91.320494
$ws.Range("I21").Value = 0.2280853681650076
$ws.Range("J21").Value = 0.2280853681650076
$ws.Range("M21").Value = 4.247957666666667
$ws.Range("N21").Value = 12.743873
$ws.Range("O21").Value = 0.1901199677306032
$ws.Range("P21").Value = 0.1901199677306032
$ws.Range("Q21").Value = 129.3085308703625
$ws.Range("R21").Value = 1163.776777833262
$ws.Range("S21").Value = 0.04336358283535399
$ws.Range("T21").Value = 0.04336358283535399

# Row 22
$ws.Range("G22").Value = 10.06688366666667
$ws.Range("H22").Value = 30.200651
$ws.Range("I22").Value = 0.07543023805979308
$ws.Range("J22").Value = 0.07543023805979308
$ws.Range("M22").Value = 1.620350333333333
$ws.Range("N22").Value = 4.861051
$ws.Range("O22").Value = 0.0725197794467048
$ws.Range("P22").Value = 0.07251977944670479
$ws.Range("Q22").Value = 16.31187830491123
$ws.Range("R22").Value = 146.806904744201
$ws.Range("S22").Value = 0.005470184227708633
$ws.Range("T22").Value = 0.005470184227708632

# Row 23
$ws.Range("G23").Value = 10.06688366666667
$ws.Range("H23").Value = 30.200651
$ws.Range("I23").Value = 0.07543023805979308
$ws.Range("J23").Value = 0.07543023805979308
$ws.Range("M23").Value = 2.170845333333333
$ws.Range("N23").Value = 6.512536
$ws.Range("O23").Value = 0.09715752300453648
$ws.Range("P23").Value = 0.09715752300453648
$ws.Range("Q23").Value = 21.85364742899289
$ws.Range("R23").Value = 196.682826860936
$ws.Range("S23").Value = 0.00732861508953201
$ws.Range("T23").Value = 0.00732861508953201

# Row 24
$ws.Range("G24").Value = 10.06688366666667
$ws.Range("H24").Value = 30.200651
$ws.Range("I24").Value = 0.07543023805979308
$ws.Range("J24").Value = 0.07543023805979308
$ws.Range("M24").Value = 13.65106133333333
$ws.Range("N24").Value = 40.953184
$ws.Range("O24").Value = 0.6109616770777183
$ws.Range("P24").Value = 0.6109616770777183
$ws.Range("Q24").Value = 137.4236463691982
$ws.Range("R24").Value = 1236.812817322784
$ws.Range("S24").Value = 0.04608498474738272
$ws.Range("T24").Value = 0.04608498474738272

# Row 25
$ws.Range("G25").Value = 10.06688366666667
$ws.Range("H25").Value = 30.200651
$ws.Range("I25").Value = 0.07543023805979308
$ws.Range("J25").Value = 0.07543023805979308
$ws.Range("M25").Value = 0.6533493333333333
$ws.Range("N25").Value = 1.960048
$ws.Range("O25").Value = 0.02924105274043717
$ws.Range("P25").Value = 0.02924105274043717
$ws.Range("Q25").Value = 6.577191732360889
$ws.Range("R25").Value = 59.194725591248
$ws.Range("S25").Value = 0.002205659569330141
$ws.Range("T25").Value = 0.002205659569330141

# Row 26
$ws.Range("G26").Value = 10.06688366666667
$ws.Range("H26").Value = 30.200651
$ws.Range("I26").Value = 0.07543023805979308
$ws.Range("J26").Value = 0.07543023805979308
$ws.Range("M26").Value = 4.247957666666667
$ws.Range("N26").Value = 12.743873
$ws.Range("O26").Value = 0.1901199677306032
$ws.Range("P26").Value = 0.1901199677306032
$ws.Range("Q26").Value = 42.76369565125812
$ws.Range("R26").Value = 384.873260861323
$ws.Range("S26").Value = 0.01434079442583958
$ws.Range("T26").Value = 0.01434079442583958
